# Updated ITA model - 2025-09-01 09:24
# Split the single "wind" fuel/commodity row in the "fuels" sheet into two
# rows: "windon" (onshore) and "windoff" (offshore).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fuels")

# Row 14 currently holds the "wind" commodity (column C), unit "twh" in
# column E. Rename it to "windon" in place.
$ws.Range("C14").Value = "windon"

# Insert a brand new row above the old row 15 (which holds the "ELC"
# commodity block) so everything below shifts down by one. The newly
# inserted row 15 becomes the "windoff" commodity, mirroring the
# structure of row 14 (commodity name in C, unit "twh" in E).
$ws.Rows.Item(15).Insert()
$ws.Range("C15").Value = "windoff"
$ws.Range("E15").Value = "twh"

# Reflect the cell the author was last focused on after the edit.
$ws.Range("E16").Select()
